# Quarterly indexing esoteric bug-fix operation
#
# The forecast-error table (rows 2-11, columns B-G) was off by one
# quarter: every row's statistics actually belonged to the quarter
# above it. This shifts all existing data rows down by one (row r
# keeps what used to be in row r-1) and fills the freed-up top row
# (row 2) with the newly computed, correct figures. Column A (the
# Q0..Q9 labels) and the header row are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B", "C", "D", "E", "F", "G")

# Snapshot current B:G values for rows 2-11 before overwriting anything.
$oldValues = @{}
for ($r = 2; $r -le 11; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value()
    }
    $oldValues[$r] = $rowVals
}

# Shift rows 3-11 down: row r now holds what used to be in row r-1.
for ($r = 11; $r -ge 3; $r--) {
    $src = $oldValues[$r - 1]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $src[$c]
    }
}

# New top row (row 2) gets the freshly computed values.
$ws.Range("B2").Value = -0.02907897629796788
$ws.Range("C2").Value = 0.3131278957257717
$ws.Range("D2").Value = 0.181524606355785
$ws.Range("E2").Value = 0.4260570458938391
$ws.Range("F2").Value = 0.43998257208981
$ws.Range("G2").Value = 15
